# Updates to second draft
#
# Applies the three connector/shape geometry tweaks from the diff
# against "Assignment 3/updatedSecondDraft.pptx" (slide 1 in this deck):
#
#   * Shape id=57  "Straight Connector 56" (a:prstGeom line)
#       off  x: 4728742 -> 4702616   (y unchanged: 2030639)
#       ext cx: 4402184 -> 4523996   (cy unchanged: 0)
#
#   * Shape id=91  "Elbow Connector 90" (bentConnector2)
#       off  x: 7093121 -> 7066995   (y unchanged: 2735313)
#       ext cx: 931275  -> 957401    (cy unchanged: 1190277)
#       + binds the connector's start point to shape id=74 ("idle"),
#         connection site index 3 (<a:stCxn id="74" idx="3"/>)
#
#   * Shape id=229 "Elbow Connector 228" (bentConnector3, rot=10800000)
#       off  x: 6470831 -> 6457770
#       off  y: 3594050 -> 3594051
#       ext cx: 1262381 -> 1249317
#       ext cy: 326597  -> 344602
#       adj1  : 1365    -> -1234

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# --- Shape id=57: "Straight Connector 56" --------------------------------
$shp57 = Get-ShapeById $s 57
$shp57.Left  = 370.28472900390625   # -> 4702616 EMU
$shp57.Width = 356.2201843261719    # -> 4523996 EMU

# --- Shape id=91: "Elbow Connector 90" ------------------------------------
$shp91 = Get-ShapeById $s 91
$shp74 = Get-ShapeById $s 74

# Bind the connector start point to the "idle" rounded-rectangle (id=74),
# connection site 3 - matches the new <a:stCxn id="74" idx="3"/>.
$shp91.ConnectorFormat.BeginConnect($shp74, 3)

$shp91.Left  = 556.4563598632812    # -> 7066995 EMU
$shp91.Width = 75.38591003417969    # -> 957401 EMU

# --- Shape id=229: "Elbow Connector 228" ----------------------------------
$shp229 = Get-ShapeById $s 229
$shp229.Left   = 508.48583984375        # -> 6457770 EMU
$shp229.Top    = 282.99615478515625     # -> 3594051 EMU
$shp229.Width  = 98.37142181396484      # -> 1249317 EMU
$shp229.Height = 27.134016036987305     # -> 344602 EMU
$shp229.Adjustments.Item(1) = -0.01234  # adj1 "val 1365" -> "val -1234"
